$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E column values (percent strings like "  +11.32%  ") are never valid numbers
# on their own (they contain "%","+","  " padding) so plain .Value assignment
# keeps them as text. D column values that look like plain decimals (e.g. "0.9986")
# would be auto-coerced to numbers by Excel, so we force those cells to Text format
# first (NumberFormat "@") -- same thing you would do by hand in Excel to keep the
# leading/trailing zeros and avoid numeric drift.

$ws.Range("D2").Value = "24.711.48"
$ws.Range("E2").Value = "  +11.32%  "

$ws.Range("D3").Value = "1.681.00"
$ws.Range("E3").Value = "  +5.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.28"
$ws.Range("E5").Value = "  +2.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9935"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3683"
$ws.Range("E7").Value = "  +1.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3415"
$ws.Range("E8").Value = "  +2.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.32"
$ws.Range("E9").Value = "  +17.10%  "

$ws.Range("E10").Value = "  +4.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07214"
$ws.Range("E11").Value = "  +4.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9951"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.093"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.14"
$ws.Range("E14").Value = "  +3.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.685"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "1.680.38"
$ws.Range("E16").Value = "  +5.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("E17").Value = "  +3.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9933"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06636"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "80.64"
$ws.Range("E20").Value = "  +6.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.36"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.052"
$ws.Range("E22").Value = "  +2.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.07"
$ws.Range("E23").Value = "  +4.36%  "

$ws.Range("D24").Value = "24.620.31"
$ws.Range("E24").Value = "  +10.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("E25").Value = "  +2.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.659"
$ws.Range("E26").Value = "  +6.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.03"
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.44"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("D29").Value = "1.866.51"
$ws.Range("E29").Value = "  +6.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.82"
$ws.Range("E30").Value = "  +4.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.223"
$ws.Range("E31").Value = "  +6.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.012"
$ws.Range("E32").Value = "  +2.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9776"
$ws.Range("E33").Value = "  +6.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08393"
$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.682"
$ws.Range("E35").Value = "  +3.47%  "

$ws.Range("E36").Value = "  +5.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06367"
$ws.Range("E37").Value = "  +5.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.289"
$ws.Range("E38").Value = "  +3.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.660"
$ws.Range("E39").Value = "  +4.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02301"
$ws.Range("E40").Value = "  +5.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.231"
$ws.Range("E41").Value = "  -0.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2080"
$ws.Range("E42").Value = "  +5.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6079"
$ws.Range("E43").Value = "  +5.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9931"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.759"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.03"
$ws.Range("E46").Value = "  +1.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5858"
$ws.Range("E47").Value = "  +5.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.60"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.000"
$ws.Range("E49").Value = "  +3.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07212"
$ws.Range("E50").Value = "  +7.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.48"
$ws.Range("E51").Value = "  +4.71%  "
